# Task2_writeup.docx edit: collapse a few "gramStart/gramEnd"-split runs back
# into single runs (proofing artifacts no longer needed), and expand the
# "How to run" instructions with an `npm install` step.

$d = $word.ActiveDocument

function Merge-ParagraphRuns {
    <#
        Re-reads a paragraph's own WordOpenXML and writes it straight back
        into the same range. Round-tripping through WordOpenXML coalesces
        adjacent same-format runs and drops now-redundant w:proofErr nodes,
        which is exactly the "remove the gramStart/gramEnd wrapped run split"
        cleanup we need - without touching any other paragraph.
    #>
    param($doc, [int]$paraIndex)

    $para = $doc.Paragraphs.Item($paraIndex)
    $range = $para.Range
    $openXml = $range.WordOpenXML

    $startTag = $openXml.IndexOf("<w:p ")
    if ($startTag -lt 0) { $startTag = $openXml.IndexOf("<w:p>") }
    $endTag = $openXml.LastIndexOf("</w:p>") + "</w:p>".Length
    $fragment = $openXml.Substring($startTag, $endTag - $startTag)

    $package = '<?xml version="1.0" standalone="yes"?>' `
        + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
        + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
        + '<pkg:xmlData>' `
        + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' `
        + '<w:body>' + $fragment + '</w:body>' `
        + '</w:document>' `
        + '</pkg:xmlData></pkg:part></pkg:package>'

    $range.InsertXML($package)
}

function Find-ParagraphIndex {
    param($doc, [string]$containsText)
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Contains($containsText)) {
            return $i
        }
    }
    return -1
}

# 1) "The data is used to compare pricing for [adjust] selling price ..."
#    -> single run, no gramStart/gramEnd around "adjust".
$idx = Find-ParagraphIndex $d "The data is used to compare pricing for"
Merge-ParagraphRuns $d $idx

# 2) "How [the] system work?" -> single run, no gramStart/gramEnd around "the".
$idx = Find-ParagraphIndex $d "system work?"
Merge-ParagraphRuns $d $idx

# 3) "[Remove duplicate] stores." -> single run, no gramStart/gramEnd.
$idx = Find-ParagraphIndex $d "Remove duplicate"
Merge-ParagraphRuns $d $idx

# 4) "... move to the file directory, then run server.js ..." gains an
#    "install node modules (npm install)" step in between.
$idx = Find-ParagraphIndex $d "move to the"
$para = $d.Paragraphs.Item($idx)
$range = $para.Range
$openXml = $range.WordOpenXML
$startTag = $openXml.IndexOf("<w:p ")
if ($startTag -lt 0) { $startTag = $openXml.IndexOf("<w:p>") }
$pprEnd = $openXml.IndexOf("</w:pPr>", $startTag) + "</w:pPr>".Length
$paraOpenAndPPr = $openXml.Substring($startTag, $pprEnd - $startTag)

$newRuns = '<w:r><w:t xml:space="preserve">On command prompt, move to the </w:t></w:r>' `
    + '<w:r><w:t>file directory,</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> install node modules ( </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>npm</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> install) </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> then run server.js (node server.js)</w:t></w:r>'

$fragment = $paraOpenAndPPr + $newRuns + '</w:p>'

$package = '<?xml version="1.0" standalone="yes"?>' `
    + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
    + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
    + '<pkg:xmlData>' `
    + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' `
    + '<w:body>' + $fragment + '</w:body>' `
    + '</w:document>' `
    + '</pkg:xmlData></pkg:part></pkg:package>'

$range.InsertXML($package)

Write-Output "Done."
